# Weekly CompStat data refresh (NYPD "cs-en-us-pbbn" report).
# Moves the report one week forward (Volume/Number + date-range header)
# and replaces the crime-statistics grid (rows 14-33) with the newly
# collected week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 31   Number  48" -> "...49" ------------------------
$hdr = $ws.Range("A8")
$numChars = $hdr.Characters(21, 2)
$numChars.Text = "49"

# --- Header: report week date range --------------------------------------
# "Report Covering the Week  11/25/2024  Through  12/1/2024"
$wk = $ws.Range("C9")
$startDate = $wk.Characters(27, 10)
$startDate.Text = "12/2/2024"
$endDate = $wk.Characters(47, 9)
$endDate.Text = "12/8/2024"

# --- Crime complaints grid (rows 14-33) -----------------------------------
$ws.Range("D14").Value = 2
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = -66.666666666666
$ws.Range("J14").Value = 66
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -13.157894736842
$ws.Range("M14").Value = -47.619047619047
$ws.Range("N14").Value = -84.722222222222
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 15
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 235
$ws.Range("J15").Value = 204
$ws.Range("K15").Value = 15.196078431372
$ws.Range("L15").Value = -0.423728813559
$ws.Range("M15").Value = 6.818181818181
$ws.Range("N15").Value = -58.916083916083
$ws.Range("C16").Value = 42
$ws.Range("E16").Value = -12.5
$ws.Range("F16").Value = 158
$ws.Range("G16").Value = 207
$ws.Range("H16").Value = -23.671497584541
$ws.Range("I16").Value = 2247
$ws.Range("J16").Value = 2376
$ws.Range("K16").Value = -5.429292929292
$ws.Range("L16").Value = -7.263722657862
$ws.Range("M16").Value = -35.39390454284
$ws.Range("N16").Value = -85.793766200923
$ws.Range("C17").Value = 68
$ws.Range("D17").Value = 66
$ws.Range("E17").Value = 3.030303030303
$ws.Range("F17").Value = 263
$ws.Range("G17").Value = 302
$ws.Range("H17").Value = -12.913907284768
$ws.Range("I17").Value = 4038
$ws.Range("J17").Value = 4005
$ws.Range("K17").Value = 0.823970037453
$ws.Range("L17").Value = 3.247251342367
$ws.Range("M17").Value = 29.589216944801
$ws.Range("N17").Value = -49.385810980195
$ws.Range("C18").Value = 26
$ws.Range("D18").Value = 42
$ws.Range("E18").Value = -38.095238095238
$ws.Range("F18").Value = 129
$ws.Range("G18").Value = 166
$ws.Range("H18").Value = -22.289156626506
$ws.Range("I18").Value = 1786
$ws.Range("J18").Value = 1929
$ws.Range("K18").Value = -7.413167444271
$ws.Range("L18").Value = -20.692717584369
$ws.Range("M18").Value = -41.307919815971
$ws.Range("N18").Value = -84.278169014084
$ws.Range("C19").Value = 85
$ws.Range("D19").Value = 101
$ws.Range("E19").Value = -15.841584158415
$ws.Range("F19").Value = 365
$ws.Range("G19").Value = 409
$ws.Range("H19").Value = -10.757946210268
$ws.Range("I19").Value = 5049
$ws.Range("J19").Value = 5481
$ws.Range("K19").Value = -7.881773399014
$ws.Range("L19").Value = -10.446967009577
$ws.Range("M19").Value = 23.296703296703
$ws.Range("N19").Value = -22.833562585969
$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 35
$ws.Range("E20").Value = -17.142857142857
$ws.Range("F20").Value = 128
$ws.Range("G20").Value = 166
$ws.Range("H20").Value = -22.89156626506
$ws.Range("I20").Value = 1652
$ws.Range("J20").Value = 1759
$ws.Range("K20").Value = -6.083001705514
$ws.Range("L20").Value = -5.761551625784
$ws.Range("M20").Value = 22.734026745913
$ws.Range("N20").Value = -81.683113427209
$ws.Range("D21").Value = 300
$ws.Range("E21").Value = -13.666666666666
$ws.Range("F21").Value = 1061
$ws.Range("G21").Value = 1271
$ws.Range("H21").Value = -16.522423288749
$ws.Range("I21").Value = 15073
$ws.Range("J21").Value = 15820
$ws.Range("K21").Value = -4.721871049304
$ws.Range("L21").Value = -7.465160537786
$ws.Range("M21").Value = -2.275674273858
$ws.Range("N21").Value = -70.857098663985
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 33.333333333333
$ws.Range("G22").Value = 29
$ws.Range("H22").Value = -34.482758620689
$ws.Range("I22").Value = 271
$ws.Range("J22").Value = 287
$ws.Range("K22").Value = -5.574912891986
$ws.Range("L22").Value = -19.345238095238
$ws.Range("M22").Value = -34.382566585956
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 24
$ws.Range("E23").Value = 29.166666666666
$ws.Range("F23").Value = 85
$ws.Range("G23").Value = 118
$ws.Range("H23").Value = -27.966101694915
$ws.Range("I23").Value = 1335
$ws.Range("J23").Value = 1479
$ws.Range("K23").Value = -9.73630831643
$ws.Range("L23").Value = -6.447091800981
$ws.Range("M23").Value = 19.516562220232
$ws.Range("C24").Value = 248
$ws.Range("D24").Value = 227
$ws.Range("E24").Value = 9.251101321585
$ws.Range("F24").Value = 927
$ws.Range("G24").Value = 882
$ws.Range("H24").Value = 5.102040816326
$ws.Range("I24").Value = 11732
$ws.Range("J24").Value = 11479
$ws.Range("K24").Value = 2.204024740831
$ws.Range("L24").Value = -7.19822812846
$ws.Range("M24").Value = 18.421318259816
$ws.Range("C25").Value = 79
$ws.Range("D25").Value = 96
$ws.Range("E25").Value = -17.708333333333
$ws.Range("F25").Value = 349
$ws.Range("G25").Value = 359
$ws.Range("H25").Value = -2.785515320334
$ws.Range("I25").Value = 5172
$ws.Range("J25").Value = 4578
$ws.Range("K25").Value = 12.975098296199
$ws.Range("L25").Value = -2.119606358819
$ws.Range("C26").Value = 109
$ws.Range("D26").Value = 97
$ws.Range("E26").Value = 12.371134020618
$ws.Range("F26").Value = 475
$ws.Range("G26").Value = 423
$ws.Range("H26").Value = 12.293144208037
$ws.Range("I26").Value = 6034
$ws.Range("J26").Value = 5730
$ws.Range("K26").Value = 5.305410122164
$ws.Range("L26").Value = 8.936631160859
$ws.Range("M26").Value = -18.766828217555
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = 12.5
$ws.Range("F27").Value = 17
$ws.Range("G27").Value = 22
$ws.Range("H27").Value = -22.727272727272
$ws.Range("I27").Value = 323
$ws.Range("J27").Value = 317
$ws.Range("K27").Value = 1.892744479495
$ws.Range("L27").Value = -9.269662921348
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = 57.142857142857
$ws.Range("F28").Value = 50
$ws.Range("G28").Value = 32
$ws.Range("H28").Value = 56.25
$ws.Range("I28").Value = 615
$ws.Range("J28").Value = 577
$ws.Range("K28").Value = 6.585788561525
$ws.Range("L28").Value = 7.7057793345
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = 16
$ws.Range("G29").Value = 24
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 229
$ws.Range("J29").Value = 230
$ws.Range("K29").Value = -0.434782608695
$ws.Range("L29").Value = -29.538461538461
$ws.Range("M29").Value = -52.390852390852
$ws.Range("N29").Value = -86.861732644865
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 66.666666666666
$ws.Range("F30").Value = 13
$ws.Range("G30").Value = 20
$ws.Range("H30").Value = -35
$ws.Range("I30").Value = 189
$ws.Range("J30").Value = 195
$ws.Range("K30").Value = -3.076923076923
$ws.Range("L30").Value = -31.021897810219
$ws.Range("M30").Value = -51.785714285714
$ws.Range("N30").Value = -87.938736439055
$ws.Range("D31").Value = 1
$ws.Range("G31").Value = 12
$ws.Range("H31").Value = -83.333333333333
$ws.Range("I31").Value = 82
$ws.Range("J31").Value = 83
$ws.Range("K31").Value = -1.204819277108
$ws.Range("L31").Value = -5.747126436781
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("F33").Value = 1
$ws.Range("H33").Value = -50
$ws.Range("J33").Value = 22
$ws.Range("K33").Value = 13.636363636363
$ws.Range("L33").Value = -10.714285714285

# Row 33 (Traffic Fatalities) gained real Week-to-Date figures this week
# (previously blank/placeholder text cells) -- give D33/E33 the same
# numeric styles already used by the other rows in this column (#,##0 and
# #,##0.0;"-"#,##0.0 respectively).
$ws.Range("D33").NumberFormat = "#,##0"
$ws.Range("E33").NumberFormat = "#,##0.0;""-""#,##0.0"
